# Edit: on slide 3 ("Сервис доставки готовой еды"), the placeholder shape
# "Объект 2" (id=3) originally holds a single paragraph:
#   "Понятия и виды доставки."
# It is changed to two paragraphs:
#   "Понятия и виды доставки."
#   "Какие есть " + "сервисы доставки"   (two separate runs)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

$target = $null
foreach ($sh in $s.Shapes) {
    if ($sh.Id -eq 3) {
        $target = $sh
    }
}

$tr = $target.TextFrame.TextRange

# Append a new paragraph after the existing one, containing the full text
# of the second line. InsertAfter() (rather than a whole-range .Text =)
# keeps the trailing endParaRPr on the new last paragraph, matching how
# PowerPoint behaves when text is typed at the end of a text box.
$tr.InsertAfter("`rКакие есть сервисы доставки")

# Re-fetch total length and split the new paragraph's text into two runs
# ("Какие есть " / "сервисы доставки") by touching the formatting of the
# second part only (re-asserting its own font size is a no-op visually,
# but forces PowerPoint to materialise it as its own <a:r> run).
$fullLen = $tr.Length
$secondRunLen = ("сервисы доставки").Length
$secondRunStart = $fullLen - $secondRunLen + 1

$secondRun = $tr.Characters($secondRunStart, $secondRunLen)
$secondRun.Font.Size = $secondRun.Font.Size
